$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.930.35"
$ws.Range("E2").Value = "  +0.47%  "

# Row 3
$ws.Range("D3").Value = "2.358.52"
$ws.Range("E3").Value = "  +1.41%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").Value = "303.27"
$ws.Range("E5").Value = "  +0.46%  "

# Row 6
$ws.Range("D6").Value = "95.11"
$ws.Range("E6").Value = "  +1.69%  "

# Row 8
$ws.Range("E8").Value = "  -0.29%  "

# Row 9
$ws.Range("E9").Value = "  -2.62%  "

# Row 10
$ws.Range("D10").Value = "34.26"
$ws.Range("E10").Value = "  +0.93%  "

# Row 11
$ws.Range("E11").Value = "  +2.01%  "

# Row 12
$ws.Range("E12").Value = "  +0.31%  "

# Row 13
$ws.Range("D13").Value = "18.47"
$ws.Range("E13").Value = "  -1.20%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "6.70"
$ws.Range("E14").Value = "  +0.13%  "

# Row 15
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.726.72"
$ws.Range("E15").Value = "  +1.40%  "

# Row 16
$ws.Range("D16").Value = "2.353.63"
$ws.Range("E16").Value = "  +1.03%  "

# Row 17
$ws.Range("E17").Value = "  +0.85%  "

# Row 18
$ws.Range("D18").Value = "42.880.08"
$ws.Range("E18").Value = "  +0.45%  "

# Row 19
$ws.Range("D19").Value = "11.95"
$ws.Range("E19").Value = "  +0.13%  "

# Row 20
$ws.Range("D20").Value = "6.25"
$ws.Range("E20").Value = "  +1.59%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0884"

# Row 22
$ws.Range("D22").Value = "67.93"
$ws.Range("E22").Value = "  +0.15%  "

# Row 23
$ws.Range("D23").Value = "234.96"
$ws.Range("E23").Value = "  -0.30%  "

# Row 24
$ws.Range("E24").Value = "  -1.90%  "

# Row 25
$ws.Range("E25").Value = "  +1.05%  "

# Row 26
$ws.Range("E26").Value = "  -0.09%  "

# Row 27
$ws.Range("D27").Value = "24.38"
$ws.Range("E27").Value = "  -0.49%  "

# Row 28
$ws.Range("D28").Value = "2.37"
$ws.Range("E28").Value = "  +6.46%  "

# Row 29
$ws.Range("D29").Value = "9.32"
$ws.Range("E29").Value = "  +2.26%  "

# Row 30
$ws.Range("D30").Value = "32.31"
$ws.Range("E30").Value = "  +3.08%  "

# Row 31
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.02%  "

# Row 32
$ws.Range("D32").Value = "5.01"
$ws.Range("E32").Value = "  +0.51%  "

# Row 33
$ws.Range("D33").Value = "17.49"
$ws.Range("E33").Value = "  -0.24%  "

# Row 34
$ws.Range("D34").Value = "0.0725"
$ws.Range("E34").Value = "  +3.79%  "

# Row 35
$ws.Range("E35").Value = "  +6.55%  "

# Row 36
$ws.Range("D36").Value = "128.65"
$ws.Range("E36").Value = "  -8.12%  "

# Row 37
$ws.Range("E37").Value = "  +0.52%  "

# Row 38
$ws.Range("D38").Value = "4.33"
$ws.Range("E38").Value = "  -0.46%  "

# Row 39
$ws.Range("D39").Value = "2.83"
$ws.Range("E39").Value = "  +3.00%  "

# Row 40
$ws.Range("D40").Value = "2.27"
$ws.Range("E40").Value = "  -2.44%  "

# Row 41
$ws.Range("E41").Value = "  -0.79%  "

# Row 42
$ws.Range("D42").Value = "20.70"
$ws.Range("E42").Value = "  -8.11%  "

# Row 43
$ws.Range("D43").Value = "1.928.37"
$ws.Range("E43").Value = "  -0.36%  "

# Row 44
$ws.Range("E44").Value = "  +0.33%  "

# Row 45
$ws.Range("E45").Value = "  +3.47%  "

# Row 46
$ws.Range("D46").Value = "9.25"
$ws.Range("E46").Value = "  -9.67%  "

# Row 47
$ws.Range("E47").Value = "  +0.09%  "

# Row 48
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "2.589.83"
$ws.Range("E48").Value = "  +1.29%  "

# Row 49
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "1.51"
$ws.Range("E49").Value = "  +3.03%  "

# Row 50
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").Value = "71.35"
$ws.Range("E50").Value = "  -1.08%  "

# Row 51
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").Value = "1.13"
$ws.Range("E51").Value = "  +1.19%  "
